# Applies the cryptos-list price/volume refresh described in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.324.93"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "2.634.12"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.24"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.42%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").Value = "2.633.67"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.361"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.76%  "
$ws.Range("D15").Value = "3.113.66"
$ws.Range("E15").Value = "  -2.66%  "
$ws.Range("E16").Value = "  -2.48%  "
$ws.Range("D17").Value = "67.203.97"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "2.627.19"
$ws.Range("E18").Value = "  -3.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "359.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.31%  "
$ws.Range("E24").Value = "  +9.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.29%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "70.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000101"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "556.01"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  -3.12%  "
$ws.Range("E35").Value = "  +4.68%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.64%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.44"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.20"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.33%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("E43").Value = "  -0.09%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").Value = "0.0₆0299"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.588"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "152.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.83"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "
